$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells to machine-readable column names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Capitalize Spanish linking particles ("de", "del", "la", "las", "el",
#    "los", "y") that appear mid-name in the state (column A) and
#    municipality (column B) text, e.g. "Pabellón de Arteaga" ->
#    "Pabellón De Arteaga". The first word of a cell is never touched.
$particles = @("de", "del", "la", "las", "el", "los", "y")

function Capitalize-Particles($s) {
    $words = $s.Split(" ")
    for ($i = 1; $i -lt $words.Length; $i++) {
        if ($particles -contains $words[$i]) {
            $words[$i] = $words[$i].Substring(0,1).ToUpper() + $words[$i].Substring(1)
        }
    }
    return [string]::Join(" ", $words)
}

$lastRow = 2350
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $textA = $cellA.Text
    if ($textA -ne "") {
        $cellA.Value = Capitalize-Particles $textA
    }
    $cellB = $ws.Cells.Item($r, 2)
    $textB = $cellB.Text
    if ($textB -ne "") {
        $cellB.Value = Capitalize-Particles $textB
    }
}

# 3. Drop the trailing footnote rows (sample size / source / author /
#    date) that lived below the data table, shrinking the used range back
#    down to A1:D2350.
$ws.Range("A2351:A2356").EntireRow.Delete() | Out-Null
